# The deck currently carries the "Integral" theme (ppt/theme/theme1.xml,
# used by the one SlideMaster) and the stock "Office Theme" colours live
# in ppt/theme/theme2.xml (only referenced by the Notes Master).
#
# The authored edit swaps the two palettes: the SlideMaster's theme
# (theme1.xml) becomes the default Office colour scheme, while the old
# Integral colours end up parked in theme2.xml. This reproduces that by
# changing the SlideMaster/Design theme's 12 colour slots to the
# standard Office palette, exactly as PowerPoint's Design gallery would
# do when a new theme/colour-scheme is applied (Master.ApplyTheme /
# ThemeColorScheme in the real object model).

function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette: built-in "Office" colour scheme.
$officeColors = @(
    (RGB 0x00 0x00 0x00), # 1  dk1
    (RGB 0xFF 0xFF 0xFF), # 2  lt1
    (RGB 0x44 0x54 0x6A), # 3  dk2
    (RGB 0xE7 0xE6 0xE6), # 4  lt2
    (RGB 0x5B 0x9B 0xD5), # 5  accent1
    (RGB 0xED 0x7D 0x31), # 6  accent2
    (RGB 0xA5 0xA5 0xA5), # 7  accent3
    (RGB 0xFF 0xC0 0x00), # 8  accent4
    (RGB 0x44 0x72 0xC4), # 9  accent5
    (RGB 0x70 0xAD 0x47), # 10 accent6
    (RGB 0x05 0x63 0xC1), # 11 hlink
    (RGB 0x95 0x4F 0x72)  # 12 folHlink
)

# The ThemeColorScheme exposes all 12 theme colour slots (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) in that order, each with a settable
# RGB property, and is backed by the presentation's one slide master
# theme part (ppt/theme/theme1.xml).
$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
